$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the cryptocurrency price/volume table (rows 2-51) to the latest
# scraped snapshot. Column D ("Price") values are stored as text in the
# source data (even when numeric-looking, e.g. "1.00" or "69.976.48" which
# isn't a valid number anyway), so a leading apostrophe is used to force
# Excel to keep them as text instead of auto-converting to a number.
# Row 2
$ws.Range("D2").Value2 = "'69.976.48"
$ws.Range("E2").Value2 = '  +2.91%  '

# Row 3
$ws.Range("D3").Value2 = "'3.405.27"
$ws.Range("E3").Value2 = '  +2.05%  '

# Row 4
$ws.Range("D4").Value2 = "'1.00"
$ws.Range("E4").Value2 = '  +0.18%  '

# Row 5
$ws.Range("D5").Value2 = "'585.27"
$ws.Range("E5").Value2 = '  +0.26%  '

# Row 6
$ws.Range("D6").Value2 = "'181.10"
$ws.Range("E6").Value2 = '  +1.88%  '

# Row 7
$ws.Range("B7").Value2 = 'XRP'
$ws.Range("C7").Value2 = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("D7").Value2 = "'0.599"
$ws.Range("E7").Value2 = '  +1.06%  '

# Row 8
$ws.Range("B8").Value2 = 'USDC'
$ws.Range("C8").Value2 = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value2 = "'1.00"
$ws.Range("E8").Value2 = '  +0.00%  '

# Row 9
$ws.Range("D9").Value2 = "'0.201"
$ws.Range("E9").Value2 = '  +9.63%  '

# Row 10
$ws.Range("D10").Value2 = "'0.595"
$ws.Range("E10").Value2 = '  +1.88%  '

# Row 11
$ws.Range("D11").Value2 = "'48.69"
$ws.Range("E11").Value2 = '  +1.75%  '

# Row 12
$ws.Range("D12").Value2 = "'0.0000288"
$ws.Range("E12").Value2 = '  +4.82%  '

# Row 13
$ws.Range("D13").Value2 = "'688.40"
$ws.Range("E13").Value2 = '  -1.42%  '

# Row 14
$ws.Range("D14").Value2 = "'8.70"
$ws.Range("E14").Value2 = '  +2.69%  '

# Row 15
$ws.Range("D15").Value2 = "'3.955.89"
$ws.Range("E15").Value2 = '  +2.13%  '

# Row 16
$ws.Range("D16").Value2 = "'70.003.13"
$ws.Range("E16").Value2 = '  +2.96%  '

# Row 17
$ws.Range("B17").Value2 = 'TRON'
$ws.Range("C17").Value2 = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").Value2 = "'0.121"
$ws.Range("E17").Value2 = '  +1.14%  '

# Row 18
$ws.Range("B18").Value2 = 'WrappedEther'
$ws.Range("C18").Value2 = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D18").Value2 = "'3.400.14"
$ws.Range("E18").Value2 = '  +1.95%  '

# Row 19
$ws.Range("D19").Value2 = "'17.77"
$ws.Range("E19").Value2 = '  +1.28%  '

# Row 20
$ws.Range("D20").Value2 = "'11.39"
$ws.Range("E20").Value2 = '  +1.97%  '

# Row 21
$ws.Range("D21").Value2 = "'0.922"
$ws.Range("E21").Value2 = '  +2.85%  '

# Row 22
$ws.Range("D22").Value2 = "'17.35"
$ws.Range("E22").Value2 = '  +2.10%  '

# Row 23
$ws.Range("D23").Value2 = "'5.37"
$ws.Range("E23").Value2 = '  -0.54%  '

# Row 24
$ws.Range("D24").Value2 = "'102.68"
$ws.Range("E24").Value2 = '  +2.15%  '

# Row 25
$ws.Range("D25").Value2 = "'3.94"
$ws.Range("E25").Value2 = '  +0.55%  '

# Row 26
$ws.Range("E26").Value2 = '  +1.05%  '

# Row 27
$ws.Range("D27").Value2 = "'9.68"
$ws.Range("E27").Value2 = '  +2.00%  '

# Row 28
$ws.Range("D28").Value2 = "'33.89"
$ws.Range("E28").Value2 = '  +2.23%  '

# Row 29
$ws.Range("D29").Value2 = "'8.86"
$ws.Range("E29").Value2 = '  +3.12%  '

# Row 30
$ws.Range("D30").Value2 = "'7.02"
$ws.Range("E30").Value2 = '  +0.52%  '

# Row 31
$ws.Range("B31").Value2 = 'dogwifhat'
$ws.Range("C31").Value2 = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D31").Value2 = "'3.71"
$ws.Range("E31").Value2 = '  +9.87%  '

# Row 32
$ws.Range("B32").Value2 = 'Cosmos'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").Value2 = "'11.15"
$ws.Range("E32").Value2 = '  +0.88%  '

# Row 33
$ws.Range("D33").Value2 = "'555.87"
$ws.Range("E33").Value2 = '  -3.54%  '

# Row 34
$ws.Range("E34").Value2 = '  +1.55%  '

# Row 35
$ws.Range("D35").Value2 = "'58.55"
$ws.Range("E35").Value2 = '  +2.15%  '

# Row 36
$ws.Range("E36").Value2 = '  +0.24%  '

# Row 37
$ws.Range("D37").Value2 = "'3.672.94"
$ws.Range("E37").Value2 = '  -2.19%  '

# Row 38
$ws.Range("D38").Value2 = "'0.141"
$ws.Range("E38").Value2 = '  +4.66%  '

# Row 39
$ws.Range("D39").Value2 = "'35.70"
$ws.Range("E39").Value2 = '  +0.44%  '

# Row 40
$ws.Range("D40").Value2 = "'0.0₃0737"
$ws.Range("E40").Value2 = '  +8.87%  '

# Row 41
$ws.Range("D41").Value2 = "'3.33"
$ws.Range("E41").Value2 = '  +4.83%  '

# Row 42
$ws.Range("D42").Value2 = "'2.72"
$ws.Range("E42").Value2 = '  +3.00%  '

# Row 43
$ws.Range("E43").Value2 = '  +4.69%  '

# Row 44
$ws.Range("D44").Value2 = "'0.340"
$ws.Range("E44").Value2 = '  +1.39%  '

# Row 45
$ws.Range("E45").Value2 = '  +1.21%  '

# Row 46
$ws.Range("D46").Value2 = "'2.68"
$ws.Range("E46").Value2 = '  +2.05%  '

# Row 47
$ws.Range("E47").Value2 = '  +0.81%  '

# Row 48
$ws.Range("B48").Value2 = 'Mantle'
$ws.Range("C48").Value2 = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").Value2 = "'1.39"
$ws.Range("E48").Value2 = '  +4.28%  '

# Row 49
$ws.Range("B49").Value2 = 'FirstDigitalUSD'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D49").Value2 = "'0.998"
$ws.Range("E49").Value2 = '  -0.31%  '

# Row 50
$ws.Range("D50").Value2 = "'130.66"
$ws.Range("E50").Value2 = '  -0.28%  '

# Row 51
$ws.Range("E51").Value2 = '  +0.30%  '

